$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on price cells whose new values would otherwise be
# auto-converted to numbers by Excel (plain decimals like "206.84").
# Percentage cells in column E and multi-dot "thousands" prices in column D
# are already unambiguous text and do not need this.
$riskyCells = @('D5','D9','D18','D22','D23','D25','D27','D30','D34','D36','D42','D45','D46','D48','D51')
foreach ($addr in $riskyCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '26.954.97'
$ws.Range('E2').Value = '  +0.18%  '
$ws.Range('D3').Value = '1.556.75'
$ws.Range('E4').Value = '  -0.16%  '
$ws.Range('D5').Value = '206.84'
$ws.Range('E5').Value = '  +0.10%  '
$ws.Range('E6').Value = '  -0.16%  '
$ws.Range('E7').Value = '  -0.16%  '
$ws.Range('E8').Value = '  +3.01%  '
$ws.Range('D9').Value = '0.247'
$ws.Range('E9').Value = '  -0.04%  '
$ws.Range('E10').Value = '  +0.88%  '
$ws.Range('E11').Value = '  +0.16%  '
$ws.Range('D12').Value = '1.778.66'
$ws.Range('E12').Value = '  +0.63%  '
$ws.Range('D13').Value = '1.556.32'
$ws.Range('E13').Value = '  +0.60%  '
$ws.Range('E14').Value = '  +1.20%  '
$ws.Range('E15').Value = '  +1.34%  '
$ws.Range('D16').Value = '26.954.09'
$ws.Range('E16').Value = '  +0.26%  '
$ws.Range('E17').Value = '  +0.53%  '
$ws.Range('D18').Value = '217.67'
$ws.Range('E18').Value = '  +1.39%  '
$ws.Range('D19').Value = '0.0₃0695'
$ws.Range('E19').Value = '  +1.73%  '
$ws.Range('E20').Value = '  +0.66%  '
$ws.Range('E21').Value = '  -0.13%  '
$ws.Range('D22').Value = '4.05'
$ws.Range('E22').Value = '  +1.04%  '
$ws.Range('D23').Value = '9.21'
$ws.Range('E23').Value = '  +0.22%  '
$ws.Range('E24').Value = '  +0.56%  '
$ws.Range('D25').Value = '153.80'
$ws.Range('E25').Value = '  +1.37%  '
$ws.Range('E26').Value = '  +0.36%  '
$ws.Range('D27').Value = '14.93'
$ws.Range('E27').Value = '  +0.56%  '
$ws.Range('E28').Value = '  +0.30%  '
$ws.Range('E29').Value = '  -0.12%  '
$ws.Range('D30').Value = '0.0468'
$ws.Range('E30').Value = '  +2.05%  '
$ws.Range('E31').Value = '  -0.52%  '
$ws.Range('E32').Value = '  -0.06%  '
$ws.Range('D33').Value = '1.423.26'
$ws.Range('E33').Value = '  +4.75%  '
$ws.Range('D34').Value = '3.09'
$ws.Range('E34').Value = '  +4.53%  '
$ws.Range('E35').Value = '  +3.29%  '
$ws.Range('D36').Value = '0.980'
$ws.Range('E36').Value = '  +2.06%  '
$ws.Range('E37').Value = '  +0.26%  '
$ws.Range('E38').Value = '  -0.14%  '
$ws.Range('E39').Value = '  -0.39%  '
$ws.Range('E40').Value = '  +0.78%  '
$ws.Range('E41').Value = '  -0.13%  '
$ws.Range('D42').Value = '5.69'
$ws.Range('E42').Value = '  +1.89%  '
$ws.Range('E43').Value = '  +3.68%  '
$ws.Range('E44').Value = '  -0.31%  '
$ws.Range('D45').Value = '64.63'
$ws.Range('E45').Value = '  +1.79%  '
$ws.Range('D46').Value = '1.75'
$ws.Range('E46').Value = '  +1.43%  '
$ws.Range('D47').Value = '1.692.26'
$ws.Range('E47').Value = '  +0.57%  '
$ws.Range('D48').Value = '87.61'
$ws.Range('E48').Value = '  +2.25%  '
$ws.Range('E49').Value = '  +1.69%  '
$ws.Range('E50').Value = '  +3.16%  '
$ws.Range('D51').Value = '0.0958'
$ws.Range('E51').Value = '  +1.13%  '
